# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated numbers (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F-column value. Note the last row (14) differs slightly
# between the two sheets in the source data.
$commonUpdates = @{
    2  = 66
    4  = 41
    6  = 2891
    8  = 1880
    10 = 95
    11 = 684
    13 = 26
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Range("F$row").Value = $commonUpdates[$row]
    }
}

# Row 14 differs between the two sheets.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F14").Value = 202

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 203
